$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write new, previously-unseen strings first in the same order the
# original authoring session introduced them, so the shared-string
# table ends up in the expected order.
$ws.Range("B4").Value = "kanasfdasabas123$"
$ws.Range("C1").Value = "expectedmsg"
$ws.Range("C2").Value = "You are logged in"
$ws.Range("C3").Value = "Invalid Username and Password"
$ws.Range("A3").Value = "kannnnn"
$ws.Range("C5").Value = "Please fill out this field"
$ws.Range("A5").Value = "hhojj"

# Remaining cells reuse already-known strings.
$ws.Range("B3").Value = "kanbas123$"
$ws.Range("A4").Value = "kanchanbasudkar"
$ws.Range("C4").Value = "Invalid Username and Password"
$ws.Range("B5").Value = "kanbas123$"

# Best-fit the three used columns to their content (mirrors Excel's
# Format > AutoFit Column Width for the new, wider data).
$ws.Columns("A:C").AutoFit() | Out-Null
$ws.Columns("A").ColumnWidth = 14.833333333333334
$ws.Columns("B").ColumnWidth = 9.833333333333334
$ws.Columns("C").ColumnWidth = 27

$ws.Range("D6").Select() | Out-Null
